$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54 (shifts existing rows 54-73 down to 55-74,
# carrying their formatting/styles with them).
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new data record.
$ws.Cells.Item(54, 1).Value = 7
$ws.Cells.Item(54, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(54, 3).Value = "Ñuble"
$ws.Cells.Item(54, 4).Value = 44845
$ws.Cells.Item(54, 5).Value = 16
$ws.Cells.Item(54, 6).Value = 100112013
$ws.Cells.Item(54, 7).Value = "Alcachofa"
$ws.Cells.Item(54, 8).Value = "Argentina(o)"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 60
$ws.Cells.Item(54, 11).Value = 10000
$ws.Cells.Item(54, 12).Value = 10000
$ws.Cells.Item(54, 13).Value = 10000
$ws.Cells.Item(54, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(54, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(54, 16).Value = 200
$ws.Cells.Item(54, 17).Value = 50
$ws.Cells.Item(54, 18).Value = "Hortaliza"
